$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows above the existing row 20. This pushes the current
# row 20 (29 / Rain Shower / 08/06/2024 / 13, all stored as text) down to
# row 23, preserving its original (text) cell types/formatting.
$ws.Rows("20:22").Insert()

# --- New row 20: 29 / Rain Shower / 08/06/2024 / 13 (numeric temp & hour) ---
$ws.Range("A20").Value = 29
$ws.Range("B20").Value = "Rain Shower"
$ws.Range("C20").NumberFormat = "@"
$ws.Range("C20").Value = "08/06/2024"
$ws.Range("C20").Style = "Normal"
$ws.Range("D20").Value = 13

# --- New row 21: 30 / Mostly Cloudy / 08/06/2024 / 14 ---
$ws.Range("A21").Value = 30
$ws.Range("B21").Value = "Mostly Cloudy"
$ws.Range("C21").NumberFormat = "@"
$ws.Range("C21").Value = "08/06/2024"
$ws.Range("C21").Style = "Normal"
$ws.Range("D21").Value = 14

# --- New row 22: 29 / Rain Shower / 08/06/2024 / 15 ---
$ws.Range("A22").Value = 29
$ws.Range("B22").Value = "Rain Shower"
$ws.Range("C22").NumberFormat = "@"
$ws.Range("C22").Value = "08/06/2024"
$ws.Range("C22").Style = "Normal"
$ws.Range("D22").Value = 15

# --- Row 23 (former row 20, shifted down): update weather + hour, keep text types ---
$ws.Range("B23").Value = "Mostly Cloudy"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "16"
$ws.Range("D23").Style = "Normal"
